$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: 01/01/1996 -> 01/01/2021 (row 8)
# Leading apostrophe forces text entry so Excel stores it as a literal
# string (matching the original cell, which also held the date as text)
# instead of auto-converting the date-looking text into a date serial.
$ws.Range("B8").Value = "'01/01/2021"
$ws.Range("C8").Value = "'01/01/2021"

# Objetivos: (row 10)
$objetivos = "- Apresentar a evolução das condições geológicas da Terra que culminaram com os recursos naturais existentes hoje, com ênfase nas reservas de combustíveis fósseis, hídricos e de minérios e a conseqüente reserva de energia advinda dessas fontes naturais. A América do Sul e do Brasil, mais especificamente, são destacados."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Programa resumido: (row 14)
$programaResumido = "- Desenvolvimento da Terra.- Recursos minerais.- Matérias-primas da grande indústria metalúrgica: metais ferrosos e não-ferrosos"
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# Programa: (row 16)
$programa = "- Desenvolvimento da Terra. - Principais Eras Geológicas. - Matérias-primas para a grande indústria metalúrgica: metais ferrosos e metais não-ferrosos."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Bibliografia: (row 22)
$bibliografia = "- MILLER, Jr. G. T. “Ciência Ambiental”,  Editora: Thomson, 2006.- ABREU, S. F. “Recursos Minerais do Brasil”, Editora: Edgard Blücher, 1973.-  SKINNER, B. J. “Recursos Minerais da Terra”, Editora: Edgard Blücher, 1996.- WICANDER, R.; MONROE, J. S. “Fundamentos de Geologia”, Editora: Cengage Learning, 2009. - PRESS, F.; Siever, R.; Jordan, T.; Grotzinger, J. “Para Entender a Terra”, Editora: Bookman,  2006.- SCHÄFER, A. “Fundamentos de Ecologia e Biogeografia de Águas Continentais”, Editora: Universidade, Porto Alegre. - Revistas especializadas e sites, dado ao caráter dinâmico das informações sobre reservas minerais e recursos naturais em geral."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
